$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 0
$ws1.Range("F3").Value = 0
$ws1.Range("F4").Value = 4802
$ws1.Range("F7").Value = 113
$ws1.Range("F12").Value = 0
$ws1.Range("F14").Value = 0
$ws1.Range("F16").Value = 0
$ws1.Range("F18").Value = 0
$ws1.Range("F19").Value = 3965
$ws1.Range("F20").Value = 0
$ws1.Range("F22").Value = 38
$ws1.Range("F23").Value = 0
$ws1.Range("F25").Value = 48
$ws1.Range("F26").Value = 0
$ws1.Range("F28").Value = 0
$ws1.Range("F29").Value = 16
$ws1.Range("F32").Value = 0
$ws1.Range("F34").Value = 0
$ws1.Range("F35").Value = 304
$ws1.Range("F36").Value = 0
$ws1.Range("F37").Value = 173
$ws1.Range("F38").Value = 0
$ws1.Range("F42").Value = 0
$ws1.Range("F44").Value = 494
$ws1.Range("F45").Value = 0
$ws1.Range("F47").Value = 0

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 227
$ws4.Range("F4").Value = 4802
$ws4.Range("F5").Value = 0
$ws4.Range("F6").Value = 155
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 93
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F14").Value = 107
$ws4.Range("F15").Value = 0
$ws4.Range("F17").Value = 0
$ws4.Range("F21").Value = 0
$ws4.Range("F23").Value = 0
$ws4.Range("F24").Value = 85
$ws4.Range("F28").Value = 0
$ws4.Range("F30").Value = 0
$ws4.Range("F32").Value = 568
$ws4.Range("F34").Value = 0
$ws4.Range("F38").Value = 173
$ws4.Range("F39").Value = 0
$ws4.Range("F41").Value = 963
$ws4.Range("F42").Value = 0
$ws4.Range("F43").Value = 0
$ws4.Range("F44").Value = 0
$ws4.Range("F45").Value = 0
$ws4.Range("F47").Value = 75
